# Generate Report for Handback
#
# The localization status report is regenerated: the c6a52559 file's
# handback has completed, so its row moves from the bottom of each sheet
# to the top (right under the header), its Status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", and its
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated. The other two rows (ffff3622... and
# ffffffcddc...) keep their data and simply shift down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Clear every hyperlink on the sheet (this engine's Range.Hyperlinks.Delete()
# clears the whole worksheet's collection) so we can re-add them, in the new
# row order, against fresh cell content.
$ov.Range("A1").Hyperlinks.Delete()

$ovRows = @(
    @{ File = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       Status = "Handed back: in sync with en-US";
       Date = "2016-03-24 07:17:57";
       Url = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md" },
    @{ File = "ffff3622fbd4-a02d-4715-8f15-e0a70f21d689.md";
       Status = "Handed back: in sync with en-US";
       Date = "2016-03-24 07:16:11";
       Url = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/ffff3622fbd4-a02d-4715-8f15-e0a70f21d689.md" },
    @{ File = "ffffffcddc5156-5980-4f20-b6b7-aa9f024a6190.md";
       Status = "Handed back: in sync with en-US";
       Date = "2016-03-24 07:16:11";
       Url = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/ffffffcddc5156-5980-4f20-b6b7-aa9f024a6190.md" }
)

$r = 2
foreach ($row in $ovRows) {
    $ov.Range("A$r").Value = $row.File
    $ov.Range("B$r").Value = $row.Status
    $ov.Range("C$r").Value = $row.Status
    $ov.Range("D$r").Value = $row.Date
    $ov.Hyperlinks.Add($ov.Range("A$r"), $row.Url, [Type]::Missing, [Type]::Missing, $row.File) | Out-Null
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A1").Hyperlinks.Delete()

$zhRows = @(
    @{ File = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       FileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       Ext = ".md";
       Status = "Handed back: in sync with en-US";
       Handoff = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.zh-cn.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8668537e49e36e4bb9dd3513a0560f755035e4d1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.zh-cn.xlf";
       HandoffDate = "2016-03-24 07:17:52";
       Target = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       TargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8668537e49e36e4bb9dd3513a0560f755035e4d1/e2e/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       Handback = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.zh-cn.xlf";
       HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8668537e49e36e4bb9dd3513a0560f755035e4d1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.zh-cn.xlf";
       HandbackDate = "2016-03-24 07:18:16";
       Reason = "Include" },
    @{ File = "ffff3622fbd4-a02d-4715-8f15-e0a70f21d689.md";
       FileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/ffff3622fbd4-a02d-4715-8f15-e0a70f21d689.md";
       Ext = ".md";
       Status = "Handed back: in sync with en-US";
       Handoff = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/141bd1ac8b3ac1645bfb2e54fdd398b5dfbbe441/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandoffDate = "2016-03-24 07:16:07";
       Target = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       TargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/73311948d399996a753690f55a2da0f5715e3438/e2e/4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       Handback = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/025c0b5478cdcf20026e85ac59807c6b801f12d0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandbackDate = "2016-03-24 07:16:30";
       Reason = "Include" },
    @{ File = "ffffffcddc5156-5980-4f20-b6b7-aa9f024a6190.md";
       FileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/ffffffcddc5156-5980-4f20-b6b7-aa9f024a6190.md";
       Ext = ".md";
       Status = "Handed back: in sync with en-US";
       Handoff = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/141bd1ac8b3ac1645bfb2e54fdd398b5dfbbe441/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandoffDate = "2016-03-24 07:16:07";
       Target = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       TargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/73311948d399996a753690f55a2da0f5715e3438/e2e/4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       Handback = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/025c0b5478cdcf20026e85ac59807c6b801f12d0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.zh-cn.xlf";
       HandbackDate = "2016-03-24 07:16:30";
       Reason = "Include" }
)

$r = 2
foreach ($row in $zhRows) {
    $zh.Range("A$r").Value = $row.File
    $zh.Range("B$r").Value = $row.Ext
    $zh.Range("C$r").Value = $row.Status
    $zh.Range("D$r").Value = $row.Handoff
    $zh.Range("E$r").Value = $row.HandoffDate
    $zh.Range("F$r").Value = $row.Target
    $zh.Range("G$r").Value = $row.Handback
    $zh.Range("H$r").Value = $row.HandbackDate
    $zh.Range("J$r").Value = $row.Reason

    $zh.Hyperlinks.Add($zh.Range("A$r"), $row.FileUrl, [Type]::Missing, [Type]::Missing, $row.File) | Out-Null
    $zh.Hyperlinks.Add($zh.Range("D$r"), $row.HandoffUrl, [Type]::Missing, [Type]::Missing, $row.Handoff) | Out-Null
    $zh.Hyperlinks.Add($zh.Range("F$r"), $row.TargetUrl, [Type]::Missing, [Type]::Missing, $row.Target) | Out-Null
    $zh.Hyperlinks.Add($zh.Range("G$r"), $row.HandbackUrl, [Type]::Missing, [Type]::Missing, $row.Handback) | Out-Null
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A1").Hyperlinks.Delete()

$deRows = @(
    @{ File = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       FileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       Ext = ".md";
       Status = "Handed back: in sync with en-US";
       Handoff = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.de-de.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e63dbf815050a694646ebf37a08fc337dfe66bc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.de-de.xlf";
       HandoffDate = "2016-03-24 07:17:57";
       Target = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       TargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0e63dbf815050a694646ebf37a08fc337dfe66bc/e2e/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.md";
       Handback = "c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.de-de.xlf";
       HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0e63dbf815050a694646ebf37a08fc337dfe66bc/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c6a52559-f5e2-425d-b0b6-ae0b4e93ee7a.98a1b9c46adbf451ed4513052adec42920e8fc69.de-de.xlf";
       HandbackDate = "2016-03-24 07:18:23";
       Reason = "Include" },
    @{ File = "ffff3622fbd4-a02d-4715-8f15-e0a70f21d689.md";
       FileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/ffff3622fbd4-a02d-4715-8f15-e0a70f21d689.md";
       Ext = ".md";
       Status = "Handed back: in sync with en-US";
       Handoff = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fe9e36be7c67940e0f34299cf711002a1252489/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandoffDate = "2016-03-24 07:16:11";
       Target = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       TargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2605eef8de5b6e4f833b328cc116a3500e46a72b/e2e/4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       Handback = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3cff72849cb5f54b6e37166698a0b537179215f2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandbackDate = "2016-03-24 07:16:37";
       Reason = "Include" },
    @{ File = "ffffffcddc5156-5980-4f20-b6b7-aa9f024a6190.md";
       FileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/15381d499dffc08e3e8b25b1a9c68fdb537f19bc/e2e/ffffffcddc5156-5980-4f20-b6b7-aa9f024a6190.md";
       Ext = ".md";
       Status = "Handed back: in sync with en-US";
       Handoff = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fe9e36be7c67940e0f34299cf711002a1252489/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandoffDate = "2016-03-24 07:16:11";
       Target = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       TargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2605eef8de5b6e4f833b328cc116a3500e46a72b/e2e/4c8092e1-ff93-42cf-b87f-8bc509240fe5.md";
       Handback = "4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3cff72849cb5f54b6e37166698a0b537179215f2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c8092e1-ff93-42cf-b87f-8bc509240fe5.e7488bd6b314a1681bd538de176a7a3450c9b7df.de-de.xlf";
       HandbackDate = "2016-03-24 07:16:37";
       Reason = "Include" }
)

$r = 2
foreach ($row in $deRows) {
    $de.Range("A$r").Value = $row.File
    $de.Range("B$r").Value = $row.Ext
    $de.Range("C$r").Value = $row.Status
    $de.Range("D$r").Value = $row.Handoff
    $de.Range("E$r").Value = $row.HandoffDate
    $de.Range("F$r").Value = $row.Target
    $de.Range("G$r").Value = $row.Handback
    $de.Range("H$r").Value = $row.HandbackDate
    $de.Range("J$r").Value = $row.Reason

    $de.Hyperlinks.Add($de.Range("A$r"), $row.FileUrl, [Type]::Missing, [Type]::Missing, $row.File) | Out-Null
    $de.Hyperlinks.Add($de.Range("D$r"), $row.HandoffUrl, [Type]::Missing, [Type]::Missing, $row.Handoff) | Out-Null
    $de.Hyperlinks.Add($de.Range("F$r"), $row.TargetUrl, [Type]::Missing, [Type]::Missing, $row.Target) | Out-Null
    $de.Hyperlinks.Add($de.Range("G$r"), $row.HandbackUrl, [Type]::Missing, [Type]::Missing, $row.Handback) | Out-Null
    $r = $r + 1
}

Write-Output "Handback report regenerated"
